$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 125, shifting the existing rows 125-131
# down to 127-133.
$ws.Rows("125:126").Insert()

# Row 125: new weekly record
$ws.Range("A125").Value = 7
$ws.Range("B125").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C125").Value = "Ñuble"
$ws.Range("D125").Value = 45041
$ws.Range("E125").Value = 16
$ws.Range("F125").Value = 100112031
$ws.Range("G125").Value = "Poroto verde"
$ws.Range("H125").Value = "Sin especificar"
$ws.Range("I125").Value = "Primera"
$ws.Range("J125").Value = 30
$ws.Range("K125").Value = 27000
$ws.Range("L125").Value = 27000
$ws.Range("M125").Value = 27000
$ws.Range("N125").Value = "`$/saco 25 kilos"
$ws.Range("O125").Value = "Región del Maule"
$ws.Range("P125").Value = 1080
$ws.Range("Q125").Value = 25
$ws.Range("R125").Value = "Hortaliza"

# Row 126: new weekly record
$ws.Range("A126").Value = 7
$ws.Range("B126").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C126").Value = "Ñuble"
$ws.Range("D126").Value = 45041
$ws.Range("E126").Value = 16
$ws.Range("F126").Value = 100112031
$ws.Range("G126").Value = "Poroto verde"
$ws.Range("H126").Value = "Sin especificar"
$ws.Range("I126").Value = "Segunda"
$ws.Range("J126").Value = 30
$ws.Range("K126").Value = 25000
$ws.Range("L126").Value = 25000
$ws.Range("M126").Value = 25000
$ws.Range("N126").Value = "`$/saco 25 kilos"
$ws.Range("O126").Value = "Región del Maule"
$ws.Range("P126").Value = 1000
$ws.Range("Q126").Value = 25
$ws.Range("R126").Value = "Hortaliza"
